$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing sheet from "Таблица умножения" to "Worksheet"
$ws.Name = "Worksheet"

# Remove the last two rows of the multiplication table (rows 5 and 6),
# shrinking the used range down to A1:E4
$ws.Rows("5:6").Delete() | Out-Null

# Keep the selection on the new bottom-right cell of the shrunk range
$ws.Range("E4").Select() | Out-Null

# Add a new, empty worksheet right after the first one
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "Worksheet 1"

# Make sure the first sheet stays the active/selected tab
$ws.Activate()
